# Update quiz results sheet:
#  - Row 2 (ID 37): ID becomes 75
#  - Row 3 (ID 38): ID becomes 76, Duration becomes 21
#  - Insert two new rows of quiz results (Deepa / I079692, Test User / I999999)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows ---
$ws.Range("A2").Value = 75
$ws.Range("A3").Value = 76
$ws.Range("G3").Value = 21

# --- Add new row 4: Deepa ---
$ws.Range("A4").Value = 77
$ws.Range("B4").Value = "Deepa"
$ws.Range("C4").Value = "I079692"
$ws.Range("D4").Value = 100
$ws.Range("E4").Value = 5
$ws.Range("F4").Value = 5
$ws.Range("G4").Value = 26
$ws.Range("H4").Value = "2025-04-27 15:47:04"

# --- Add new row 5: Test User ---
$ws.Range("A5").Value = 78
$ws.Range("B5").Value = "Test User"
$ws.Range("C5").Value = "I999999"
$ws.Range("D5").Value = 80
$ws.Range("E5").Value = 5
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = "2025-04-27 15:43:16"
